$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Data Warehouse Engineer"
$ws.Range("B3").Value = "GO-JEK"
$ws.Range("C3").Value = "https://id.indeed.com//cmp/Pt.-Go--jek-Indonesia-2"
$ws.Range("A4").Value = "Lead Instructors - Le Wagon Data Science Bootcamp"
$ws.Range("B4").Value = "Le Wagon Bali"
$ws.Range("C4").Value = "Link is not available"
$ws.Range("A5").Value = "ShopeePay Backend Engineer [Experienced]"
$ws.Range("B5").Value = "Shopee"
$ws.Range("C5").Value = "https://id.indeed.com//cmp/Shopee"
$ws.Range("A6").Value = "Back End Developer"
$ws.Range("B6").Value = "PT Generasi Teknologi Buana"
$ws.Range("C6").Value = "Link is not available"
$ws.Range("A8").Value = "Associate Backend Engineer"
$ws.Range("B8").Value = "KeDA Tech"
$ws.Range("C8").Value = "Link is not available"
$ws.Range("A9").Value = "QA Engineer"
$ws.Range("B9").Value = "StyleTheory"
$ws.Range("C9").Value = "https://id.indeed.com//cmp/Styletheory"
$ws.Range("A10").Value = "Developer / Programmer"
$ws.Range("B10").Value = "StrategArt"
$ws.Range("A11").Value = "Back End Developer - Kompas.id"
$ws.Range("B11").Value = "Kompas Gramedia"
$ws.Range("C11").Value = "https://id.indeed.com//cmp/Kompas-Gramedia-Group"
$ws.Range("A12").Value = "Back End Developer"
$ws.Range("B12").Value = "JULO"
$ws.Range("A13").Value = "Python Junior Programmer"
$ws.Range("B13").Value = "Sonar Social Media Monitoring Platform"
$ws.Range("C13").Value = "Link is not available"
$ws.Range("A14").Value = "Problem Generator Developer"
$ws.Range("B14").Value = "Zenius Education"
$ws.Range("A15").Value = "Data Warehouse Engineer - GoPay"
$ws.Range("B15").Value = "GO-JEK"
$ws.Range("C15").Value = "https://id.indeed.com//cmp/Pt.-Go--jek-Indonesia-2"
$ws.Range("A16").Value = "Web Developer"
$ws.Range("B16").Value = "Great Giant Foods"
$ws.Range("A17").Value = "Application Developer"
$ws.Range("B17").Value = "Bank Mega"
$ws.Range("C17").Value = "https://id.indeed.com//cmp/Bank-Mega"
$ws.Range("A18").Value = "Senior Machine Learning"
$ws.Range("B18").Value = "Alodokter"
$ws.Range("C18").Value = "https://id.indeed.com//cmp/Alodokter-2"
$ws.Range("A19").Value = "Product Support Engineer"
$ws.Range("B19").Value = "Awan Tunai"
$ws.Range("A20").Value = "Software Engineer - Data Platform"
$ws.Range("B20").Value = "Cermati.com"
$ws.Range("C20").Value = "https://id.indeed.com//cmp/PT-Dwi-Cermat-Indonesia-1"
$ws.Range("A21").Value = "Project Manager / Jr. Project Manager / ERP Senior Consultan..."
$ws.Range("B21").Value = "HashMicro"
$ws.Range("A22").Value = "Remote Senior Web Engineer"
$ws.Range("B22").Value = "Scopic"
$ws.Range("C22").Value = "Link is not available"
$ws.Range("A23").Value = "HRIS Developer"
$ws.Range("B23").Value = "Binabusana Internusa"
$ws.Range("C23").Value = "Link is not available"
$ws.Range("A24").Value = "Financial Service SRE Engineer [Entry Level]"
$ws.Range("A25").Value = "Web Developer"
$ws.Range("B25").Value = "1rstWAP"
$ws.Range("C25").Value = "Link is not available"
$ws.Range("A26").Value = "Full Stack Developer"
$ws.Range("B26").Value = "Global Talentlytica"
$ws.Range("C26").Value = "Link is not available"
$ws.Range("A27").Value = "Python Programmer"
$ws.Range("B27").Value = "1rstWAP"
$ws.Range("A28").Value = "Technical Solution Developer"
$ws.Range("B28").Value = "PT Mastersystem Infotama"
$ws.Range("C28").Value = "https://id.indeed.com//cmp/Pt.-Mastersystem-Infotama"
$ws.Range("A29").Value = "Backend Developer"
$ws.Range("B29").Value = "Pintek ID"
$ws.Range("A30").Value = "Test Engineer"
$ws.Range("B30").Value = "Quipper"
$ws.Range("C30").Value = "Link is not available"
$ws.Range("A31").Value = "Engineering and Technology - Site Reliability Engineer"
$ws.Range("B31").Value = "Shopee"
$ws.Range("C31").Value = "https://id.indeed.com//cmp/Shopee"
$ws.Range("A32").Value = "Machine Learning Engineer"
$ws.Range("B32").Value = "Nomura Research Institute Indonesia"
$ws.Range("A33").Value = "IT WEB DEVELOPER"
$ws.Range("B33").Value = "Ismaya Group"
$ws.Range("C33").Value = "https://id.indeed.com//cmp/Ismaya-Group-1"
$ws.Range("A34").Value = "Back End Developer"
$ws.Range("B34").Value = "Akseleran"
$ws.Range("C34").Value = "Link is not available"
$ws.Range("A35").Value = "Back End Developer"
$ws.Range("B35").Value = "Renos.id"
$ws.Range("A36").Value = "ShopeePay Backend Engineer [Leader]"
$ws.Range("B36").Value = "Shopee"
$ws.Range("C36").Value = "https://id.indeed.com//cmp/Shopee"
$ws.Range("A37").Value = "Web Developer"
$ws.Range("B37").Value = "Binabusana Internusa"
$ws.Range("A38").Value = "IT Programmer Analyst"
$ws.Range("B38").Value = "PT BSR Indonesia"
$ws.Range("C38").Value = "Link is not available"
$ws.Range("A39").Value = "Golang Developer (Back End)"
$ws.Range("B39").Value = "PT Lunaria Annua Teknologi (KoinWorks)"
$ws.Range("A40").Value = "Business Intelligence Developer E-Commerce"
$ws.Range("B40").Value = "Kompas Gramedia"
$ws.Range("C40").Value = "https://id.indeed.com//cmp/Kompas-Gramedia-Group"
$ws.Range("A41").Value = "IT Developer"
$ws.Range("B41").Value = "MNC"
$ws.Range("A42").Value = "Engineering and Technology - Back End Engineer, Payment Proc..."
$ws.Range("B42").Value = "Shopee"
$ws.Range("C42").Value = "https://id.indeed.com//cmp/Shopee"
$ws.Range("A43").Value = "Senior Frontend Developer"
$ws.Range("B43").Value = "Ensoft"
$ws.Range("C43").Value = "Link is not available"
$ws.Range("A44").Value = "Full Stack Developer"
$ws.Range("B44").Value = "IndoSterling Technomedia"
$ws.Range("C44").Value = "Link is not available"
$ws.Range("A45").Value = "Test Engineer"
$ws.Range("B45").Value = "LINE Plus corporation"
$ws.Range("C45").Value = "Link is not available"
$ws.Range("A46").Value = "Senior Data Warehouse Engineer"
$ws.Range("B46").Value = "GO-JEK"
$ws.Range("C46").Value = "https://id.indeed.com//cmp/Pt.-Go--jek-Indonesia-2"
$ws.Range("A47").Value = "Unity Developer"
$ws.Range("B47").Value = "Alegrium"
$ws.Range("C47").Value = "Link is not available"
$ws.Range("A48").Value = "Test Engineer"
$ws.Range("B48").Value = "LINE Plus corporation"
$ws.Range("C48").Value = "Link is not available"
$ws.Range("A49").Value = "Full-Stack Developer"
$ws.Range("B49").Value = "DDTC"
$ws.Range("A50").Value = "Software Quality Assurance - Manual Testing"
$ws.Range("B50").Value = "Cermati.com"
$ws.Range("C50").Value = "https://id.indeed.com//cmp/PT-Dwi-Cermat-Indonesia-1"
$ws.Range("A51").Value = "Senior Data Warehouse Engineer"
$ws.Range("B51").Value = "GO-JEK"
$ws.Range("C51").Value = "https://id.indeed.com//cmp/Pt.-Go--jek-Indonesia-2"
$ws.Range("A52").Value = "Fullstack Developer - Javascript, Python, Golang, NodeJS, Re..."
$ws.Range("B52").Value = "Michael Page"
$ws.Range("C52").Value = "https://id.indeed.com//cmp/Michael-Page"
$ws.Range("A53").Value = "Customer Solutions Consultant, Infrastructure Modernization,..."
$ws.Range("B53").Value = "Google"
$ws.Range("C53").Value = "https://id.indeed.com//cmp/Google"
$ws.Range("A54").Value = "IT Engineering Manager"
$ws.Range("B54").Value = "Ajaib"
$ws.Range("A55").Value = "System Administrator"
$ws.Range("B55").Value = "Jawasoft"
$ws.Range("A56").Value = "IT Developer"
$ws.Range("B56").Value = "MNC"
$ws.Range("C56").Value = "Link is not available"
$ws.Range("A57").Value = "Engineer: Software Developer"
$ws.Range("B57").Value = "NTT Ltd"
$ws.Range("C57").Value = "Link is not available"
$ws.Range("A58").Value = "Back End Developer"
$ws.Range("B58").Value = "Renos.id"
$ws.Range("A59").Value = "Software Developer"
$ws.Range("B59").Value = "Terrindo Bumi Raya"
$ws.Range("C59").Value = "Link is not available"
$ws.Range("A60").Value = "Full Stack Developer"
$ws.Range("B60").Value = "Global Talentlytica"
$ws.Range("C60").Value = "Link is not available"
$ws.Range("A61").Value = "Senior Integration Developer"
$ws.Range("B61").Value = "GO-JEK"
$ws.Range("C61").Value = "https://id.indeed.com//cmp/Pt.-Go--jek-Indonesia-2"
$ws.Range("A62").Value = "Software Quality Assurance - Manual Testing"
$ws.Range("B62").Value = "Cermati.com"
$ws.Range("C62").Value = "https://id.indeed.com//cmp/PT-Dwi-Cermat-Indonesia-1"
$ws.Range("A63").Value = "Technical Operations Engineer"
$ws.Range("B63").Value = "byOrange"
$ws.Range("C63").Value = "Link is not available"
$ws.Range("A64").Value = "Frontend Engineer (Freshgraduate, Senior, Principal, & Senio..."
$ws.Range("B64").Value = "Cermati.com"
$ws.Range("C64").Value = "https://id.indeed.com//cmp/PT-Dwi-Cermat-Indonesia-1"
$ws.Range("A65").Value = "(Singapore Corp) Python Software Developer"
$ws.Range("B65").Value = "MatchaTalent"
$ws.Range("A66").Value = "IT Engineering Manager"
$ws.Range("B66").Value = "Ajaib"
$ws.Range("A67").Value = "System Administrator"
$ws.Range("B67").Value = "Jawasoft"
$ws.Range("A68").Value = "Python Developer"
$ws.Range("A69").Value = "Financial Service SRE Engineer [Experienced]"
$ws.Range("B69").Value = "Shopee"
$ws.Range("C69").Value = "https://id.indeed.com//cmp/Shopee"
$ws.Range("A70").Value = "Software Engineer"
$ws.Range("B70").Value = "Alterra"
$ws.Range("A71").Value = "Finance MIS Analyst - Financial Services"
$ws.Range("B71").Value = "GO-JEK"
$ws.Range("C71").Value = "https://id.indeed.com//cmp/Pt.-Go--jek-Indonesia-2"
$ws.Range("A72").Value = "Back End Engineer"
$ws.Range("B72").Value = "Ruangguru"
$ws.Range("C72").Value = "https://id.indeed.com//cmp/PT-Ruang-Raya-Indonesia-(ruangguru)"
$ws.Range("A73").Value = "Server-side Engineer"
$ws.Range("B73").Value = "LINE Plus corporation"
$ws.Range("A74").Value = "Unity Developer"
$ws.Range("B74").Value = "Alegrium"
$ws.Range("C74").Value = "Link is not available"
$ws.Range("A75").Value = "Engineering and Technology - System Quality Assurance"
$ws.Range("B75").Value = "Shopee"
$ws.Range("C75").Value = "https://id.indeed.com//cmp/Shopee"
$ws.Range("A76").Value = "IT Engineer"
$ws.Range("B76").Value = "PT Bank Central Asia Tbk"
$ws.Range("C76").Value = "https://id.indeed.com//cmp/PT-Bank-Central-Asia-Tbk"
